$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.847.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.104.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.31%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.103.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.36%  "
$ws.Range("E9").Value = "  -5.00%  "
$ws.Range("E10").Value = "  -8.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.14%  "
$ws.Range("E12").Value = "  -6.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -10.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.614.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.37%  "
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.849.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.102.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.22%  "
$ws.Range("E22").Value = "  -8.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.00%  "
$ws.Range("E24").Value = "  -8.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.03%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.83%  "
$ws.Range("E29").Value = "  -12.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.115"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.96%  "
$ws.Range("E31").Value = "  -5.38%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.09%  "
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("E36").Value = "  -9.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0742"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "461.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -14.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0392"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.26%  "
$ws.Range("E42").Value = "  -8.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.840.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.91%  "
$ws.Range("E45").Value = "  -10.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -13.60%  "
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.79%  "
$ws.Range("E50").Value = "  -5.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.69%  "
